$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column E ("Result") for all data rows now reads "N/A" instead of "PASS"/"FAIL"
$ws.Range("E2:E7").Value = "N/A"

# Move the active selection to B4 (was E8)
$ws.Range("B4").Select()
